{"js": "// Replace every arithmetic expression in the 20x5 results table with its\n// updated expression. The new values are taken from the target-revision\n// OOXML (each <w:t> inside the table, read in row-major document order);\n// one \"old\" expression (\"42+17=59\") repeats for two different cells, so a\n// global find/replace would be ambiguous -- we instead overwrite the whole\n// table grid positionally via `Table.values`, which preserves each cell's\n// existing paragraph/run formatting (fonts, size, alignment, etc.) and only\n// swaps the run text.\n\nconst newValues = [\n  [\"63-53=10\", \"38+6=44\", \"88-35=53\", \"66-57=9\", \"47+39=86\"],\n  [\"93-88=5\", \"36+17=53\", \"91-77=14\", \"37+42=79\", \"87+10=97\"],\n  [\"83-49=34\", \"89-26=63\", \"4+3=7\", \"55+3=58\", \"58+36=94\"],\n  [\"0+46=46\", \"1+71=72\", \"9+48=57\", \"99-91=8\", \"32+56=88\"],\n  [\"38-37=1\", \"42+1=43\", \"11+17=28\", \"80-57=23\", \"26+46=72\"],\n  [\"37+60=97\", \"91-36=55\", \"34+60=94\", \"9+22=31\", \"45+34=79\"],\n  [\"9+51=60\", \"0+91=91\", \"62-16=46\", \"60-16=44\", \"62-30=32\"],\n  [\"9+9=18\", \"76+9=85\", \"2+58=60\", \"32-5=27\", \"21-4=17\"],\n  [\"90-60=30\", \"20+13=33\", \"11+55=66\", \"24+68=92\", \"68+13=81\"],\n  [\"11+2=13\", \"9+27=36\", \"28+37=65\", \"52+2=54\", \"44-21=23\"],\n  [\"96-13=83\", \"99-78=21\", \"26+19=45\", \"94-1=93\", \"51+0=51\"],\n  [\"17-6=11\", \"87-40=47\", \"9+37=46\", \"58-54=4\", \"36+21=57\"],\n  [\"63+27=90\", \"50+14=64\", \"70+26=96\", \"32+51=83\", \"11+13=24\"],\n  [\"64-49=15\", \"93-23=70\", \"50-6=44\", \"76-21=55\", \"19-12=7\"],\n  [\"70+26=96\", \"0+27=27\", \"64+31=95\", \"25-18=7\", \"95-20=75\"],\n  [\"86-83=3\", \"39-16=23\", \"36-9=27\", \"62-60=2\", \"19+1=20\"],\n  [\"5+76=81\", \"1+81=82\", \"20-16=4\", \"80+10=90\", \"82-65=17\"],\n  [\"69-7=62\", \"19+37=56\", \"73-51=22\", \"39-37=2\", \"11+23=34\"],\n  [\"4+88=92\", \"87-49=38\", \"4+84=88\", \"39+58=97\", \"7+55=62\"],\n  [\"58+0=58\", \"22-10=12\", \"52-6=46\", \"32+64=96\", \"97-42=55\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount, values\");\nawait context.sync();\n\nif (table.rowCount !== newValues.length) {\n  throw new Error(\n    `Expected ${newValues.length} rows, found ${table.rowCount}.`\n  );\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace each arithmetic expression in the 20x5 results table with the new\n# expression, in document (row-major) order. Index-based (not global find/\n# replace) so the one duplicated \"old\" text (\"42+17=59\", cells 22 and 23)\n# still maps to two different new values.\n$d = $word.ActiveDocument\n\n$newValues = @(\n    \"63-53=10\",\n    \"38+6=44\",\n    \"88-35=53\",\n    \"66-57=9\",\n    \"47+39=86\",\n    \"93-88=5\",\n    \"36+17=53\",\n    \"91-77=14\",\n    \"37+42=79\",\n    \"87+10=97\",\n    \"83-49=34\",\n    \"89-26=63\",\n    \"4+3=7\",\n    \"55+3=58\",\n    \"58+36=94\",\n    \"0+46=46\",\n    \"1+71=72\",\n    \"9+48=57\",\n    \"99-91=8\",\n    \"32+56=88\",\n    \"38-37=1\",\n    \"42+1=43\",\n    \"11+17=28\",\n    \"80-57=23\",\n    \"26+46=72\",\n    \"37+60=97\",\n    \"91-36=55\",\n    \"34+60=94\",\n    \"9+22=31\",\n    \"45+34=79\",\n    \"9+51=60\",\n    \"0+91=91\",\n    \"62-16=46\",\n    \"60-16=44\",\n    \"62-30=32\",\n    \"9+9=18\",\n    \"76+9=85\",\n    \"2+58=60\",\n    \"32-5=27\",\n    \"21-4=17\",\n    \"90-60=30\",\n    \"20+13=33\",\n    \"11+55=66\",\n    \"24+68=92\",\n    \"68+13=81\",\n    \"11+2=13\",\n    \"9+27=36\",\n    \"28+37=65\",\n    \"52+2=54\",\n    \"44-21=23\",\n    \"96-13=83\",\n    \"99-78=21\",\n    \"26+19=45\",\n    \"94-1=93\",\n    \"51+0=51\",\n    \"17-6=11\",\n    \"87-40=47\",\n    \"9+37=46\",\n    \"58-54=4\",\n    \"36+21=57\",\n    \"63+27=90\",\n    \"50+14=64\",\n    \"70+26=96\",\n    \"32+51=83\",\n    \"11+13=24\",\n    \"64-49=15\",\n    \"93-23=70\",\n    \"50-6=44\",\n    \"76-21=55\",\n    \"19-12=7\",\n    \"70+26=96\",\n    \"0+27=27\",\n    \"64+31=95\",\n    \"25-18=7\",\n    \"95-20=75\",\n    \"86-83=3\",\n    \"39-16=23\",\n    \"36-9=27\",\n    \"62-60=2\",\n    \"19+1=20\",\n    \"5+76=81\",\n    \"1+81=82\",\n    \"20-16=4\",\n    \"80+10=90\",\n    \"82-65=17\",\n    \"69-7=62\",\n    \"19+37=56\",\n    \"73-51=22\",\n    \"39-37=2\",\n    \"11+23=34\",\n    \"4+88=92\",\n    \"87-49=38\",\n    \"4+84=88\",\n    \"39+58=97\",\n    \"7+55=62\",\n    \"58+0=58\",\n    \"22-10=12\",\n    \"52-6=46\",\n    \"32+64=96\",\n    \"97-42=55\"\n)\n\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\nif (($rows * $cols) -ne $newValues.Count) {\n    throw \"Expected a $($newValues.Count)-cell table, found $rows x $cols.\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i]\n        $i++\n    }\n}\n\n"}
